$d = $word.ActiveDocument

# Update the date heading paragraph (direct range assignment keeps it scoped
# to this paragraph only, and preserves the run/paragraph formatting).
$d.Paragraphs(1).Range.Text = "2025-08-31 Sunday"

# Update the division problems in the table, cell by cell. Several of the new
# values duplicate another cell's OLD value (e.g. cell(1,1)'s new text is the
# same string as cell(5,1)'s old text), so a document-wide Find/Replace would
# clobber the wrong cell. Assigning each Cell's Range.Text directly keeps every
# edit scoped to exactly that cell.
$t = $d.Tables.Item(1)
$t.Cell(1,1).Range.Text = "45÷9=5, 0"
$t.Cell(1,2).Range.Text = "53÷4=13, 1"
$t.Cell(1,3).Range.Text = "63÷6=10, 3"
$t.Cell(1,4).Range.Text = "62÷9=6, 8"
$t.Cell(1,5).Range.Text = "46÷8=5, 6"
$t.Cell(5,1).Range.Text = "50÷7=7, 1"
$t.Cell(5,2).Range.Text = "89÷5=17, 4"
$t.Cell(5,3).Range.Text = "17÷4=4, 1"
$t.Cell(5,4).Range.Text = "74÷2=37, 0"
$t.Cell(5,5).Range.Text = "48÷5=9, 3"
$t.Cell(9,1).Range.Text = "19÷5=3, 4"
$t.Cell(9,2).Range.Text = "62÷2=31, 0"
$t.Cell(9,3).Range.Text = "65÷2=32, 1"
$t.Cell(9,4).Range.Text = "97÷9=10, 7"
$t.Cell(9,5).Range.Text = "80÷5=16, 0"
$t.Cell(13,1).Range.Text = "43÷4=10, 3"
$t.Cell(13,2).Range.Text = "56÷9=6, 2"
$t.Cell(13,3).Range.Text = "39÷9=4, 3"
$t.Cell(13,4).Range.Text = "25÷2=12, 1"
$t.Cell(13,5).Range.Text = "82÷3=27, 1"
$t.Cell(17,1).Range.Text = "94÷4=23, 2"
$t.Cell(17,2).Range.Text = "36÷6=6, 0"
$t.Cell(17,3).Range.Text = "10÷7=1, 3"
$t.Cell(17,4).Range.Text = "15÷4=3, 3"
$t.Cell(17,5).Range.Text = "44÷7=6, 2"
